# backup before deadline cleaning
# Refresh the recalculated result values in row 2 of Sheet1 (res_sum_Q_outrB_RA).
# These are the updated output numbers after re-running the underlying model;
# most differ only in floating-point-noise digits, a handful moved materially.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "E2"  = 1.1900000000000002
    "G2"  = 1.5089999999999999
    "H2"  = 1.913
    "M2"  = 1.8310000000000002
    "N2"  = 1.7429999999999999
    "R2"  = 1.3539999999999999
    "W2"  = 1.7509999999999999
    "Y2"  = 0.86899999999999999
    "Z2"  = 1.675
    "AA2" = 2.294
    "AB2" = 2.0939999999999999
    "AC2" = 2.0879999999999996
    "AE2" = 1.7410000000000001
    "AG2" = 2.9720000000000004
    "AK2" = 2.14
    "AN2" = 1.0920000000000001
    "AQ2" = 0.05
    "AS2" = 1.5839999999999999
    "AV2" = 0.84899999999999998
    "AW2" = 0.91700000000000004
    "BB2" = 0.85300000000000009
    "BC2" = 1.228
    "BD2" = 0.82400000000000007
    "BE2" = 1.6700000000000002
    "BI2" = 2.0420000000000003
    "BJ2" = 2.7960000000000003
    "BM2" = 2.0220000000000002
    "BQ2" = 2.2490000000000001
    "BS2" = 2.0019999999999998
    "BU2" = 1.7100000000000002
    "BZ2" = 2.4560000000000004
    "CF2" = 2.8089999999999997
    "CK2" = 2.4329999999999998
    "CQ2" = 1.167
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
